# petty-cashBook-2021.xlsx — 30-Jan-2021 midday update
# Adds new ledger entries for 29-Jan and 30-Jan (Buku KAS HARIAN "Sheet1"
# tab), extending existing transactions on rows 29-31 and filling rows
# 33-47 which were previously blank placeholder rows (only carrying the
# running-balance formula in column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- extend existing 28-Jan-2021 entries (rows 29-31) ---------------------
$ws.Range("D29").Formula = "=60000+280000"
$ws.Range("C30").Formula = "=13320000+18450000+40274000+29120000+16368000+14266000+500000+20616000"
$ws.Range("D31").Formula = "=29370000+14266000+2000000+301000"

# --- row 33: CHEQUE RECEIVED ----------------------------------------------
$ws.Range("B33").Value = "CHEQUE RECEIVED"
$ws.Range("D33").Formula = "=2162000"

# --- row 34: PRIVE - andreas -----------------------------------------------
$ws.Range("B34").Value = "PRIVE - andreas"
$ws.Range("D34").Value = 5000000

# --- row 35: SALES - cash/retail -------------------------------------------
$ws.Range("B35").Value = "SALES - cash/retail"
$ws.Range("C35").Formula = "=111913025-79834025-20616000"

# --- row 36: SERVICE - pintu ------------------------------------------------
$ws.Range("B36").Value = "SERVICE - pintu"
$ws.Range("D36").Value = 350000

# --- row 37: SELISIH - lebih -------------------------------------------------
$ws.Range("B37").Value = "SELISIH - lebih"
$ws.Range("C37").Value = 35000

# --- row 38: SETOR KE BANK ---------------------------------------------------
$ws.Range("B38").Value = "SETOR KE BANK"
$ws.Range("D38").Value = 111000000

# --- row 39: new day, 29-Jan-2021 (44225), Wages Expense --------------------
$ws.Range("A39").Value = 44225
$ws.Range("B39").Value = "Wages Expense"
$ws.Range("D39").Formula = "=60000+300000+300000"

# --- row 40: A/R -------------------------------------------------------------
$ws.Range("B40").Value = "A/R"
$ws.Range("C40").Formula = "=5000000+17850000+6490000+26626500"

# --- row 41: TRANSFER BCA -----------------------------------------------------
$ws.Range("B41").Value = "TRANSFER BCA"
$ws.Range("D41").Formula = "=17850000+180000+430000+6490000+425000+2249500"

# --- row 42: FREIGHT IN --------------------------------------------------------
$ws.Range("B42").Value = "FREIGHT IN"
$ws.Range("D42").Formula = "=2715500"

# --- row 43: SALES - cash/retail -----------------------------------------------
$ws.Range("B43").Value = "SALES - cash/retail"
$ws.Range("C43").Formula = "=31726525+1041975-26626500"

# --- row 44: SELISIH - lebih -----------------------------------------------------
$ws.Range("B44").Value = "SELISIH - lebih"
$ws.Range("C44").Value = 83500

# --- row 45: SETOR KE BANK ---------------------------------------------------------
$ws.Range("B45").Value = "SETOR KE BANK"
$ws.Range("D45").Value = 20000000

# --- row 46: new day, 30-Jan-2021 (44226), Wages Expense ----------------------------
$ws.Range("A46").Value = 44226
$ws.Range("B46").Value = "Wages Expense"
$ws.Range("D46").Formula = "=60000"

# --- row 47: BELI ban -----------------------------------------------------------------
$ws.Range("B47").Value = "BELI ban"
$ws.Range("D47").Formula = "=698500"

# --- view state: scroll the frozen pane and move the active selection ------------------
$ws.Range("A45").Select()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A66").Select()
